$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "57.990.20"
$ws.Range("E2").Value = "  -0.37%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.450.52"
$ws.Range("E3").Value = "  -1.04%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "511.78"
$ws.Range("E5").Value = "  -1.71%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "130.20"
$ws.Range("E6").Value = "  -1.19%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  -0.16%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.552"
$ws.Range("E8").Value = "  -1.40%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.468.27"
$ws.Range("E9").Value = "  -0.46%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.157"
$ws.Range("E10").Value = "  +0.16%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0961"
$ws.Range("E11").Value = "  -3.48%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.20"
$ws.Range("E12").Value = "  -2.93%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.331"
$ws.Range("E13").Value = "  -3.69%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.885.06"
$ws.Range("E14").Value = "  -1.07%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "57.952.27"
$ws.Range("E15").Value = "  -0.30%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "22.02"
$ws.Range("E16").Value = "  -0.41%  "

$ws.Range("E17").Value = "  -2.32%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.460.15"
$ws.Range("E18").Value = "  -0.70%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.54"
$ws.Range("E19").Value = "  -2.80%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "318.70"
$ws.Range("E20").Value = "  -0.76%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.13"
$ws.Range("E21").Value = "  -1.45%  "

$ws.Range("E22").Value = "  +0.17%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.96"
$ws.Range("E23").Value = "  +3.33%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.30"
$ws.Range("E24").Value = "  -1.30%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.405"
$ws.Range("E25").Value = "  -1.29%  "

$ws.Range("E26").Value = "  -0.67%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.160"
$ws.Range("E27").Value = "  -1.17%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.27"
$ws.Range("E28").Value = "  -1.89%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "168.06"
$ws.Range("E29").Value = "  +1.00%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0731"
$ws.Range("E30").Value = "  -3.20%  "

$ws.Range("B31").Value = "Fetch.AI"
$ws.Range("C31").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.18"
$ws.Range("E31").Value = "  -1.42%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.22"
$ws.Range("E32").Value = "  -2.06%  "

$ws.Range("B33").Value = "PancakeSwap"
$ws.Range("C33").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.67"
$ws.Range("E33").Value = "  -2.22%  "

$ws.Range("E34").Value = "  -0.07%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.994"
$ws.Range("E35").Value = "  -0.37%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.82"
$ws.Range("E36").Value = "  -1.93%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.28"
$ws.Range("E37").Value = "  -3.96%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.93"
$ws.Range("E38").Value = "  -1.89%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.59"
$ws.Range("E39").Value = "  -0.03%  "

$ws.Range("E40").Value = "  -1.67%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.763"
$ws.Range("E41").Value = "  -4.14%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "271.96"
$ws.Range("E42").Value = "  -1.80%  "

$ws.Range("E43").Value = "  -0.88%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.39"
$ws.Range("E44").Value = "  -2.84%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.588"
$ws.Range("E45").Value = "  -1.49%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0914"
$ws.Range("E46").Value = "  +0.59%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "119.93"
$ws.Range("E47").Value = "  -5.19%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0490"
$ws.Range("E48").Value = "  -0.33%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "17.30"
$ws.Range("E49").Value = "  -3.46%  "

$ws.Range("E50").Value = "  -2.22%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "16.72"
$ws.Range("E51").Value = "  -2.88%  "
